{"js": "// Diff adds, right before the sectPr (i.e. at the very end of the body):\n//   1) a new paragraph with two runs: \"4\" (rFonts hint=\"eastAsia\") and \"56\"\n//   2) a new, completely empty paragraph\n//\n// Build that exact OOXML fragment and insert it at the end of the body so\n// Word splices it in as brand-new paragraphs (rather than merging the runs\n// into the existing last paragraph).\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:hint=\"eastAsia\"/>\n              </w:rPr>\n              <w:t>4</w:t>\n            </w:r>\n            <w:r>\n              <w:t>56</w:t>\n            </w:r>\n          </w:p>\n          <w:p/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst endOfBody = context.document.body.getRange(\"End\");\nendOfBody.insertOoxml(ooxml, \"End\");\nawait context.sync();\n", "ps1": "# Diff adds, right before the sectPr (i.e. at the very end of the body):\n#   1) a new paragraph with two runs: \"4\" (rFonts hint=\"eastAsia\") and \"56\"\n#   2) a new, completely empty paragraph\n#\n# Build that exact OOXML fragment and insert it at a zero-length range\n# collapsed to the end of the document's content, so Word splices it in as\n# brand-new paragraphs (rather than merging the runs into the existing last\n# paragraph's run).\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Collapse(0)  # wdCollapseEnd\n\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:hint=\"eastAsia\"/>\n              </w:rPr>\n              <w:t>4</w:t>\n            </w:r>\n            <w:r>\n              <w:t>56</w:t>\n            </w:r>\n          </w:p>\n          <w:p/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($ooxml)\n"}
